$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ContactUs")
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "ToDoItems"
$newSheet.Range("A1").Value = "items"
$newSheet.Range("A2").Value = "1. new item 1"
$newSheet.Range("A3").Value = "new item 2"
$newSheet.Range("A4").Value = "new &&!*@"
$newSheet.Range("A5").Value = 123
$newSheet.Range("A6").Select() | Out-Null
